# pom4.xml MF Scenario Correction
# Refresh the "CreatedSuites" / "UpdatedSuites" test-run log sheets with a
# newer batch of timestamped suite names, trimming the stale rows and
# re-styling the first logged "Suite-" entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CreatedSuites sheet
# ---------------------------------------------------------------------
$created = $wb.Worksheets.Item("CreatedSuites")

# A2 becomes the newest "Suite-" timestamp and picks up a bigger, custom
# font (sz 12, SF-Pro-Display-Regular) - which also grows that row's height.
$created.Range("A2").Value = "Suite-08:45:12"
$created.Range("A2").Font.Size = 12
$created.Range("A2").Font.Name = "SF-Pro-Display-Regular"
$created.Rows(2).RowHeight = 15.5

# Refresh the remaining three kept rows with newer timestamps ...
$created.Range("A3").Value = "Suite-09:00:10"
$created.Range("A4").Value = "Royal-09:01:11"
$created.Range("A5").Value = "Royal-09:02:03"

# ... and drop the trailing three rows that no longer apply.
$created.Rows("6:8").Delete()

# ---------------------------------------------------------------------
# UpdatedSuites sheet
# ---------------------------------------------------------------------
$updated = $wb.Worksheets.Item("UpdatedSuites")

$updated.Range("A2").Value = "Suite-09:00:10"
$updated.Range("B2").Value = "Royal-09:01:11"
$updated.Range("A3").Value = "Royal-09:01:11"
$updated.Range("B3").Value = "Royal-09:02:03"

# Drop the trailing two rows that no longer apply.
$updated.Rows("4:5").Delete()

# Restore the sheet's active selection.
$updated.Activate()
$updated.Range("B13").Select()
